$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - keeps the existing bold/border style that was already on A1
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "Как вас зовут?"
$ws.Range("C1").Value = "Какой ваш любимый цвет?"
$ws.Range("D1").Value = "Какой ваш любимый фильм?"
$ws.Range("E1").Value = "Какой ваш любимый вид спорта?"

# Copy the header style from A1 (bold + border) onto B1:E1
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore values after paste-special in case formats overwrote them (paste special with formats only should not touch values)
$ws.Range("B1").Value = "Как вас зовут?"
$ws.Range("C1").Value = "Какой ваш любимый цвет?"
$ws.Range("D1").Value = "Какой ваш любимый фильм?"
$ws.Range("E1").Value = "Какой ваш любимый вид спорта?"

# Row 2
$ws.Range("A2").Value = "LowIQMulti"
$ws.Range("B2").Value = "denis"
$ws.Range("C2").Value = "Красный"
$ws.Range("D2").Value = "UIIUII"
$ws.Range("E2").Value = "Баскетбол"

# Row 3
$ws.Range("A3").Value = "LowIQMulti"
$ws.Range("B3").Value = "ty"
$ws.Range("C3").Value = "Красный"
$ws.Range("D3").Value = "ty"
$ws.Range("E3").Value = "Баскетбол"

# Row 4
$ws.Range("A4").Value = "drus1k0"
$ws.Range("B4").Value = "Bob"
$ws.Range("C4").Value = "Синий"
$ws.Range("D4").Value = "Hohlo"
$ws.Range("E4").Value = "Футбол"

# Row 5
$ws.Range("A5").Value = "LowIQMulti"
$ws.Range("B5").Value = "Егор"
$ws.Range("C5").Value = "Красный"
$ws.Range("D5").Value = "неет"
$ws.Range("E5").Value = "Баскетбол"

# Row 6
$ws.Range("A6").Value = "GasBillt"
$ws.Range("B6").Value = "Dis"
$ws.Range("C6").Value = "Синий"
$ws.Range("D6").Value = "Tg"
$ws.Range("E6").Value = "Футбол"
